$d = $word.ActiveDocument

# 1. "A Telegram Bot available via @sectionche" -> "...@sectionchec"
$d.Content.Find.Execute("A Telegram Bot available via @sectionche", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A Telegram Bot available via @sectionchec", 2)

# 2. "cker_bot in Telegram Messaging App" -> "ker_bot in Telegram Messaging App"
$d.Content.Find.Execute("cker_bot in Telegram Messaging App", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ker_bot in Telegram Messaging App", 2)

# 3. "more than 500 people" -> "more than 1000 people"
$d.Content.Find.Execute("more than 500 people", $true, $false, $false, $false, $false,
                         $true, 1, $false, "more than 1000 people", 2)

# 4. Nudge a decorative underline image's vertical offset by -1 EMU
#    (wp:posOffset 7827293 -> 7827292), a layout recalculation artifact.
#    1 EMU = 1/914400 inch = 1/12700 point, so set Top with full precision.
$shape = $d.Shapes.Item(8)
$shape.Top = 7827292 / 914400.0 * 72.0
